# Update countries & provincias Spain
# Applies the daily COVID-19 data refresh:
#  - bumps the "last updated" timestamp in A1
#  - updates case counters for several countries (rows identified by
#    their position in the country table)
#  - Rumania's totals overtook Suecia/Belgica, so those three rows are
#    re-sorted (Rumania moves up, Suecia/Belgica shift down one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 12:26"

# --- Estados Unidos (row 4) ---------------------------------------------
$ws.Range("B4").Value = 6048404
$ws.Range("C4").Value = 1770
$ws.Range("D4").Value = 3348744
$ws.Range("E4").Value = 2514826
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 184834

# --- Iran (row 15) --------------------------------------------------------
$ws.Range("B15").Value = 369911
$ws.Range("C15").Value = 2115
$ws.Range("D15").Value = 318270
$ws.Range("E15").Value = 30392
$ws.Range("G15").Value = 112
$ws.Range("H15").Value = 21249

# --- Banglades (row 18) --------------------------------------------------
$ws.Range("B18").Value = 306794
$ws.Range("C18").Value = 2211
$ws.Range("D18").Value = 196836
$ws.Range("E18").Value = 105784
$ws.Range("G18").Value = 47
$ws.Range("H18").Value = 4174

# --- Rumania / Suecia / Belgica re-sort (rows 39-41) ---------------------
# Rumania's updated total (84468) now exceeds Suecia (83898) and Belgica
# (83500), so it moves from row 41 up to row 39; Suecia and Belgica each
# shift down one row, keeping the table sorted by total cases descending.
$ws.Range("A39").Value = "Rumania"
$ws.Range("B39").Value = 84468
$ws.Range("C39").Value = 1318
$ws.Range("D39").Value = 37056
$ws.Range("E39").Value = 43905
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 48
$ws.Range("H39").Value = 3507

$ws.Range("A40").Value = "Suecia"
$ws.Range("B40").Value = 83898
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 5820

$ws.Range("A41").Value = "Belgica"
$ws.Range("B41").Value = 83500
$ws.Range("C41").Value = 470
$ws.Range("D41").Value = 18360
$ws.Range("E41").Value = 55256
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 9884

# --- Consejo Danes para los Refugiados (row 91) ---------------------------
$ws.Range("B91").Value = 9994
$ws.Range("C91").Value = 79
$ws.Range("D91").Value = 9035
$ws.Range("E91").Value = 704

# --- Finlandia (row 101) --------------------------------------------------
$ws.Range("B101").Value = 8042
$ws.Range("C101").Value = 23
$ws.Range("E101").Value = 507
